$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.926.86"
$ws.Range("E2").Value = "  -0.31%  "

$ws.Range("D3").Value = "2.788.24"
$ws.Range("E3").Value = "  -2.07%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "357.92"
$ws.Range("E5").Value = "  -0.80%  "

$ws.Range("D6").Value = "109.66"
$ws.Range("E6").Value = "  -3.00%  "

$ws.Range("E7").Value = "  -1.97%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  -2.60%  "

$ws.Range("D10").Value = "40.38"
$ws.Range("E10").Value = "  -2.82%  "

$ws.Range("E11").Value = "  +1.61%  "

$ws.Range("E12").Value = "  -1.88%  "

$ws.Range("D13").Value = "19.51"
$ws.Range("E13").Value = "  -3.78%  "

$ws.Range("E14").Value = "  -3.22%  "

$ws.Range("D15").Value = "3.228.87"
$ws.Range("E15").Value = "  -2.07%  "

$ws.Range("D16").Value = "2.790.41"
$ws.Range("E16").Value = "  -1.95%  "

$ws.Range("D17").Value = "0.950"
$ws.Range("E17").Value = "  +1.86%  "

$ws.Range("D18").Value = "51.880.41"
$ws.Range("E18").Value = "  -0.29%  "

$ws.Range("D19").Value = "7.45"
$ws.Range("E19").Value = "  -2.35%  "

$ws.Range("E20").Value = "  -2.51%  "

$ws.Range("D21").Value = "13.17"
$ws.Range("E21").Value = "  -3.05%  "

$ws.Range("D22").Value = "0.0₃0977"
$ws.Range("E22").Value = "  -2.22%  "

$ws.Range("D23").Value = "270.93"
$ws.Range("E23").Value = "  +0.40%  "

$ws.Range("D24").Value = "70.20"
$ws.Range("E24").Value = "  -0.38%  "

$ws.Range("E25").Value = "  -3.78%  "

$ws.Range("D26").Value = "26.48"
$ws.Range("E26").Value = "  -2.66%  "

$ws.Range("E27").Value = "  -0.07%  "

$ws.Range("E28").Value = "  +17.71%  "

$ws.Range("D29").Value = "10.32"
$ws.Range("E29").Value = "  -0.78%  "

$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "2.14"
$ws.Range("E30").Value = "  -5.13%  "

$ws.Range("B31").Value = "VeChain"
$ws.Range("C31").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D31").Value = "0.0469"
$ws.Range("E31").Value = "  -0.99%  "

$ws.Range("D32").Value = "52.01"

$ws.Range("D33").Value = "34.83"
$ws.Range("E33").Value = "  -2.52%  "

$ws.Range("E34").Value = "  -3.57%  "

$ws.Range("D35").Value = "0.0846"
$ws.Range("E35").Value = "  -0.27%  "

$ws.Range("D36").Value = "5.21"
$ws.Range("E36").Value = "  -5.94%  "

$ws.Range("E37").Value = "  -0.06%  "

$ws.Range("D38").Value = "18.77"
$ws.Range("E38").Value = "  +0.60%  "

$ws.Range("E39").Value = "  -3.46%  "

$ws.Range("E40").Value = "  -4.60%  "

$ws.Range("D41").Value = "2.59"
$ws.Range("E41").Value = "  +1.38%  "

$ws.Range("E42").Value = "  -2.05%  "

$ws.Range("E43").Value = "  -1.83%  "

$ws.Range("D44").Value = "119.70"
$ws.Range("E44").Value = "  -4.42%  "

$ws.Range("D45").Value = "21.94"
$ws.Range("E45").Value = "  -7.09%  "

$ws.Range("D46").Value = "2.079.55"
$ws.Range("E46").Value = "  -1.59%  "

$ws.Range("D47").Value = "3.28"
$ws.Range("E47").Value = "  -4.51%  "

$ws.Range("E48").Value = "  -1.99%  "

$ws.Range("D49").Value = "5.77"
$ws.Range("E49").Value = "  -4.49%  "

$ws.Range("D50").Value = "0.951"
$ws.Range("E50").Value = "  -3.37%  "

$ws.Range("E51").Value = "  -4.53%  "
